$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shift existing D:K data to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Populate new columns D and E with new quarterly financial data (Q4 2018 / Q3 2018)
$ws.Range("D7:E7").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D7:E7").Font.Name = "Verdana"
$ws.Range("D7:E7").Font.Size = 12
$ws.Range("D7:E7").Font.Bold = $true
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373

$ws.Range("D8:E8").NumberFormat = "#,##0"
$ws.Range("D8:E8").Font.Name = "Verdana"
$ws.Range("D8:E8").Font.Size = 12
$ws.Range("D8:E8").Font.Bold = $false
$ws.Range("D8:E8").HorizontalAlignment = -4152
$ws.Range("D8").Value = 1460700
$ws.Range("E8").Value = 1409000

$ws.Range("D9:E9").NumberFormat = "#,##0"
$ws.Range("D9:E9").Font.Name = "Verdana"
$ws.Range("D9:E9").Font.Size = 12
$ws.Range("D9:E9").Font.Bold = $false
$ws.Range("D9:E9").HorizontalAlignment = -4152
$ws.Range("D9").Value = 254100
$ws.Range("E9").Value = 262000

$ws.Range("D10:E10").NumberFormat = "#,##0"
$ws.Range("D10:E10").Font.Name = "Verdana"
$ws.Range("D10:E10").Font.Size = 12
$ws.Range("D10:E10").Font.Bold = $false
$ws.Range("D10:E10").HorizontalAlignment = -4152
$ws.Range("D10").Value = 1206600
$ws.Range("E10").Value = 1147000

$ws.Range("D12:E12").NumberFormat = "#,##0"
$ws.Range("D12:E12").Font.Name = "Verdana"
$ws.Range("D12:E12").Font.Size = 12
$ws.Range("D12:E12").Font.Bold = $false
$ws.Range("D12:E12").HorizontalAlignment = -4152
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"

$ws.Range("D13:E13").NumberFormat = "#,##0"
$ws.Range("D13:E13").Font.Name = "Verdana"
$ws.Range("D13:E13").Font.Size = 12
$ws.Range("D13:E13").Font.Bold = $false
$ws.Range("D13:E13").HorizontalAlignment = -4152
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

$ws.Range("D14:E14").NumberFormat = "#,##0"
$ws.Range("D14:E14").Font.Name = "Verdana"
$ws.Range("D14:E14").Font.Size = 12
$ws.Range("D14:E14").Font.Bold = $false
$ws.Range("D14:E14").HorizontalAlignment = -4152
$ws.Range("D14").Value = -143900
$ws.Range("E14").Value = 0

$ws.Range("D15:E15").NumberFormat = "#,##0"
$ws.Range("D15:E15").Font.Name = "Verdana"
$ws.Range("D15:E15").Font.Size = 12
$ws.Range("D15:E15").Font.Bold = $false
$ws.Range("D15:E15").HorizontalAlignment = -4152
$ws.Range("D15").Value = 329100
$ws.Range("E15").Value = 316200

$ws.Range("D17:E17").NumberFormat = "#,##0"
$ws.Range("D17:E17").Font.Name = "Verdana"
$ws.Range("D17:E17").Font.Size = 12
$ws.Range("D17:E17").Font.Bold = $false
$ws.Range("D17:E17").HorizontalAlignment = -4152
$ws.Range("D17").Value = 562800
$ws.Range("E17").Value = 691600

$ws.Range("D18:E18").NumberFormat = "#,##0"
$ws.Range("D18:E18").Font.Name = "Verdana"
$ws.Range("D18:E18").Font.Size = 12
$ws.Range("D18:E18").Font.Bold = $false
$ws.Range("D18:E18").HorizontalAlignment = -4152
$ws.Range("D18").Value = 897900
$ws.Range("E18").Value = 717400

$ws.Range("D20:E20").NumberFormat = "#,##0"
$ws.Range("D20:E20").Font.Name = "Verdana"
$ws.Range("D20:E20").Font.Size = 12
$ws.Range("D20:E20").Font.Bold = $false
$ws.Range("D20:E20").HorizontalAlignment = -4152
$ws.Range("D20").Value = 139600
$ws.Range("E20").Value = 124300

$ws.Range("D21:E21").NumberFormat = "#,##0"
$ws.Range("D21:E21").Font.Name = "Verdana"
$ws.Range("D21:E21").Font.Size = 12
$ws.Range("D21:E21").Font.Bold = $false
$ws.Range("D21:E21").HorizontalAlignment = -4152
$ws.Range("D21").Value = 1377600
$ws.Range("E21").Value = 1177300

$ws.Range("D22:E22").NumberFormat = "#,##0"
$ws.Range("D22:E22").Font.Name = "Verdana"
$ws.Range("D22:E22").Font.Size = 12
$ws.Range("D22:E22").Font.Bold = $false
$ws.Range("D22:E22").HorizontalAlignment = -4152
$ws.Range("D22").Value = 204300
$ws.Range("E22").Value = 199500

$ws.Range("D23:E23").NumberFormat = "#,##0"
$ws.Range("D23:E23").Font.Name = "Verdana"
$ws.Range("D23:E23").Font.Size = 12
$ws.Range("D23:E23").Font.Bold = $false
$ws.Range("D23:E23").HorizontalAlignment = -4152
$ws.Range("D23").Value = 833200
$ws.Range("E23").Value = 642200

$ws.Range("D24:E24").NumberFormat = "#,##0"
$ws.Range("D24:E24").Font.Name = "Verdana"
$ws.Range("D24:E24").Font.Size = 12
$ws.Range("D24:E24").Font.Bold = $false
$ws.Range("D24:E24").HorizontalAlignment = -4152
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0

$ws.Range("D25:E25").NumberFormat = "#,##0"
$ws.Range("D25:E25").Font.Name = "Verdana"
$ws.Range("D25:E25").Font.Size = 12
$ws.Range("D25:E25").Font.Bold = $false
$ws.Range("D25:E25").HorizontalAlignment = -4152
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

$ws.Range("D26:E26").NumberFormat = "#,##0"
$ws.Range("D26:E26").Font.Name = "Verdana"
$ws.Range("D26:E26").Font.Size = 12
$ws.Range("D26:E26").Font.Bold = $false
$ws.Range("D26:E26").HorizontalAlignment = -4152
$ws.Range("D26").Value = 833200
$ws.Range("E26").Value = 642200

$ws.Range("D27:E27").NumberFormat = "#,##0"
$ws.Range("D27:E27").Font.Name = "Verdana"
$ws.Range("D27:E27").Font.Size = 12
$ws.Range("D27:E27").Font.Bold = $false
$ws.Range("D27:E27").HorizontalAlignment = -4152
$ws.Range("D27").Value = 712800
$ws.Range("E27").Value = 556300

$ws.Range("D28:E28").NumberFormat = "#,##0"
$ws.Range("D28:E28").Font.Name = "Verdana"
$ws.Range("D28:E28").Font.Size = 12
$ws.Range("D28:E28").Font.Bold = $false
$ws.Range("D28:E28").HorizontalAlignment = -4152
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

$ws.Range("D29:E29").NumberFormat = "#,##0"
$ws.Range("D29:E29").Font.Name = "Verdana"
$ws.Range("D29:E29").Font.Size = 12
$ws.Range("D29:E29").Font.Bold = $false
$ws.Range("D29:E29").HorizontalAlignment = -4152
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

$ws.Range("D30:E30").NumberFormat = "#,##0"
$ws.Range("D30:E30").Font.Name = "Verdana"
$ws.Range("D30:E30").Font.Size = 12
$ws.Range("D30:E30").Font.Bold = $false
$ws.Range("D30:E30").HorizontalAlignment = -4152
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

$ws.Range("D31:E31").NumberFormat = "#,##0"
$ws.Range("D31:E31").Font.Name = "Verdana"
$ws.Range("D31:E31").Font.Size = 12
$ws.Range("D31:E31").Font.Bold = $false
$ws.Range("D31:E31").HorizontalAlignment = -4152
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

$ws.Range("D32:E32").NumberFormat = "#,##0"
$ws.Range("D32:E32").Font.Name = "Verdana"
$ws.Range("D32:E32").Font.Size = 12
$ws.Range("D32:E32").Font.Bold = $false
$ws.Range("D32:E32").HorizontalAlignment = -4152
$ws.Range("D32").Value = -139600
$ws.Range("E32").Value = -124300

$ws.Range("D33:E33").NumberFormat = "#,##0"
$ws.Range("D33:E33").Font.Name = "Verdana"
$ws.Range("D33:E33").Font.Size = 12
$ws.Range("D33:E33").Font.Bold = $false
$ws.Range("D33:E33").HorizontalAlignment = -4152
$ws.Range("D33").Value = 712800
$ws.Range("E33").Value = 556300

$ws.Range("D34:E34").NumberFormat = "#,##0"
$ws.Range("D34:E34").Font.Name = "Verdana"
$ws.Range("D34:E34").Font.Size = 12
$ws.Range("D34:E34").Font.Bold = $false
$ws.Range("D34:E34").HorizontalAlignment = -4152
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

$ws.Range("D35:E35").NumberFormat = "#,##0"
$ws.Range("D35:E35").Font.Name = "Verdana"
$ws.Range("D35:E35").Font.Size = 12
$ws.Range("D35:E35").Font.Bold = $false
$ws.Range("D35:E35").HorizontalAlignment = -4152
$ws.Range("D35").Value = 712800
$ws.Range("E35").Value = 556300

$ws.Range("D38:E38").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D38:E38").Font.Name = "Verdana"
$ws.Range("D38:E38").Font.Size = 12
$ws.Range("D38:E38").Font.Bold = $true
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373

$ws.Range("D41:E41").NumberFormat = "#,##0"
$ws.Range("D41:E41").Font.Name = "Verdana"
$ws.Range("D41:E41").Font.Size = 12
$ws.Range("D41:E41").Font.Bold = $false
$ws.Range("D41:E41").HorizontalAlignment = -4152
$ws.Range("D41").Value = 514300
$ws.Range("E41").Value = 695700

$ws.Range("D42:E42").NumberFormat = "#,##0"
$ws.Range("D42:E42").Font.Name = "Verdana"
$ws.Range("D42:E42").Font.Size = 12
$ws.Range("D42:E42").Font.Bold = $false
$ws.Range("D42:E42").HorizontalAlignment = -4152
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0

$ws.Range("D43:E43").NumberFormat = "#,##0"
$ws.Range("D43:E43").Font.Name = "Verdana"
$ws.Range("D43:E43").Font.Size = 12
$ws.Range("D43:E43").Font.Bold = $false
$ws.Range("D43:E43").HorizontalAlignment = -4152
$ws.Range("D43").Value = 763800
$ws.Range("E43").Value = 722700

$ws.Range("D44:E44").NumberFormat = "#,##0"
$ws.Range("D44:E44").Font.Name = "Verdana"
$ws.Range("D44:E44").Font.Size = 12
$ws.Range("D44:E44").Font.Bold = $false
$ws.Range("D44:E44").HorizontalAlignment = -4152
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0

$ws.Range("D45:E45").NumberFormat = "#,##0"
$ws.Range("D45:E45").Font.Name = "Verdana"
$ws.Range("D45:E45").Font.Size = 12
$ws.Range("D45:E45").Font.Bold = $false
$ws.Range("D45:E45").HorizontalAlignment = -4152
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0

$ws.Range("D46:E46").NumberFormat = "#,##0"
$ws.Range("D46:E46").Font.Name = "Verdana"
$ws.Range("D46:E46").Font.Size = 12
$ws.Range("D46:E46").Font.Bold = $false
$ws.Range("D46:E46").HorizontalAlignment = -4152
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0

$ws.Range("D47:E47").NumberFormat = "#,##0"
$ws.Range("D47:E47").Font.Name = "Verdana"
$ws.Range("D47:E47").Font.Size = 12
$ws.Range("D47:E47").Font.Bold = $false
$ws.Range("D47:E47").HorizontalAlignment = -4152
$ws.Range("D47").Value = 3989900
$ws.Range("E47").Value = 4058300

$ws.Range("D48:E48").NumberFormat = "#,##0"
$ws.Range("D48:E48").Font.Name = "Verdana"
$ws.Range("D48:E48").Font.Size = 12
$ws.Range("D48:E48").Font.Bold = $false
$ws.Range("D48:E48").HorizontalAlignment = -4152
$ws.Range("D48").Value = 24208100
$ws.Range("E48").Value = 24304900

$ws.Range("D49:E49").NumberFormat = "#,##0"
$ws.Range("D49:E49").Font.Name = "Verdana"
$ws.Range("D49:E49").Font.Size = 12
$ws.Range("D49:E49").Font.Bold = $false
$ws.Range("D49:E49").HorizontalAlignment = -4152
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0

$ws.Range("D50:E50").NumberFormat = "#,##0"
$ws.Range("D50:E50").Font.Name = "Verdana"
$ws.Range("D50:E50").Font.Size = 12
$ws.Range("D50:E50").Font.Bold = $false
$ws.Range("D50:E50").HorizontalAlignment = -4152
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

$ws.Range("D51:E51").NumberFormat = "#,##0"
$ws.Range("D51:E51").Font.Name = "Verdana"
$ws.Range("D51:E51").Font.Size = 12
$ws.Range("D51:E51").Font.Bold = $false
$ws.Range("D51:E51").HorizontalAlignment = -4152
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

$ws.Range("D52:E52").NumberFormat = "#,##0"
$ws.Range("D52:E52").Font.Name = "Verdana"
$ws.Range("D52:E52").Font.Size = 12
$ws.Range("D52:E52").Font.Bold = $false
$ws.Range("D52:E52").HorizontalAlignment = -4152
$ws.Range("D52").Value = 1210000
$ws.Range("E52").Value = 1298000

$ws.Range("D53:E53").NumberFormat = "#,##0"
$ws.Range("D53:E53").Font.Name = "Verdana"
$ws.Range("D53:E53").Font.Size = 12
$ws.Range("D53:E53").Font.Bold = $false
$ws.Range("D53:E53").HorizontalAlignment = -4152
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

$ws.Range("D54:E54").NumberFormat = "#,##0"
$ws.Range("D54:E54").Font.Name = "Verdana"
$ws.Range("D54:E54").Font.Size = 12
$ws.Range("D54:E54").Font.Bold = $false
$ws.Range("D54:E54").HorizontalAlignment = -4152
$ws.Range("D54").Value = 30686200
$ws.Range("E54").Value = 31079700

$ws.Range("D57:E57").NumberFormat = "#,##0"
$ws.Range("D57:E57").Font.Name = "Verdana"
$ws.Range("D57:E57").Font.Size = 12
$ws.Range("D57:E57").Font.Bold = $false
$ws.Range("D57:E57").HorizontalAlignment = -4152
$ws.Range("D57").Value = 1316900
$ws.Range("E57").Value = 1268100

$ws.Range("D58:E58").NumberFormat = "#,##0"
$ws.Range("D58:E58").Font.Name = "Verdana"
$ws.Range("D58:E58").Font.Size = 12
$ws.Range("D58:E58").Font.Bold = $false
$ws.Range("D58:E58").HorizontalAlignment = -4152
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0

$ws.Range("D59:E59").NumberFormat = "#,##0"
$ws.Range("D59:E59").Font.Name = "Verdana"
$ws.Range("D59:E59").Font.Size = 12
$ws.Range("D59:E59").Font.Bold = $false
$ws.Range("D59:E59").HorizontalAlignment = -4152
$ws.Range("D59").Value = 1536100
$ws.Range("E59").Value = 1534600

$ws.Range("D60:E60").NumberFormat = "#,##0"
$ws.Range("D60:E60").Font.Name = "Verdana"
$ws.Range("D60:E60").Font.Size = 12
$ws.Range("D60:E60").Font.Bold = $false
$ws.Range("D60:E60").HorizontalAlignment = -4152
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0

$ws.Range("D61:E61").NumberFormat = "#,##0"
$ws.Range("D61:E61").Font.Name = "Verdana"
$ws.Range("D61:E61").Font.Size = 12
$ws.Range("D61:E61").Font.Bold = $false
$ws.Range("D61:E61").HorizontalAlignment = -4152
$ws.Range("D61").Value = 23305500
$ws.Range("E61").Value = 23678300

$ws.Range("D62:E62").NumberFormat = "#,##0"
$ws.Range("D62:E62").Font.Name = "Verdana"
$ws.Range("D62:E62").Font.Size = 12
$ws.Range("D62:E62").Font.Bold = $false
$ws.Range("D62:E62").HorizontalAlignment = -4152
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0

$ws.Range("D63:E63").NumberFormat = "#,##0"
$ws.Range("D63:E63").Font.Name = "Verdana"
$ws.Range("D63:E63").Font.Size = 12
$ws.Range("D63:E63").Font.Bold = $false
$ws.Range("D63:E63").HorizontalAlignment = -4152
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

$ws.Range("D64:E64").NumberFormat = "#,##0"
$ws.Range("D64:E64").Font.Name = "Verdana"
$ws.Range("D64:E64").Font.Size = 12
$ws.Range("D64:E64").Font.Bold = $false
$ws.Range("D64:E64").HorizontalAlignment = -4152
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

$ws.Range("D65:E65").NumberFormat = "#,##0"
$ws.Range("D65:E65").Font.Name = "Verdana"
$ws.Range("D65:E65").Font.Size = 12
$ws.Range("D65:E65").Font.Bold = $false
$ws.Range("D65:E65").HorizontalAlignment = -4152
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

$ws.Range("D66:E66").NumberFormat = "#,##0"
$ws.Range("D66:E66").Font.Name = "Verdana"
$ws.Range("D66:E66").Font.Size = 12
$ws.Range("D66:E66").Font.Bold = $false
$ws.Range("D66:E66").HorizontalAlignment = -4152
$ws.Range("D66").Value = 27389500
$ws.Range("E66").Value = 27696100

$ws.Range("D68:E68").NumberFormat = "#,##0"
$ws.Range("D68:E68").Font.Name = "Verdana"
$ws.Range("D68:E68").Font.Size = 12
$ws.Range("D68:E68").Font.Bold = $false
$ws.Range("D68:E68").HorizontalAlignment = -4152
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

$ws.Range("D69:E69").NumberFormat = "#,##0"
$ws.Range("D69:E69").Font.Name = "Verdana"
$ws.Range("D69:E69").Font.Size = 12
$ws.Range("D69:E69").Font.Bold = $false
$ws.Range("D69:E69").HorizontalAlignment = -4152
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

$ws.Range("D70:E70").NumberFormat = "#,##0"
$ws.Range("D70:E70").Font.Name = "Verdana"
$ws.Range("D70:E70").Font.Size = 12
$ws.Range("D70:E70").Font.Bold = $false
$ws.Range("D70:E70").HorizontalAlignment = -4152
$ws.Range("D70").Value = 42700
$ws.Range("E70").Value = 42800

$ws.Range("D71:E71").NumberFormat = "#,##0"
$ws.Range("D71:E71").Font.Name = "Verdana"
$ws.Range("D71:E71").Font.Size = 12
$ws.Range("D71:E71").Font.Bold = $false
$ws.Range("D71:E71").HorizontalAlignment = -4152
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

$ws.Range("D72:E72").NumberFormat = "#,##0"
$ws.Range("D72:E72").Font.Name = "Verdana"
$ws.Range("D72:E72").Font.Size = 12
$ws.Range("D72:E72").Font.Bold = $false
$ws.Range("D72:E72").HorizontalAlignment = -4152
$ws.Range("D72").Value = -4893100
$ws.Range("E72").Value = -4896800

$ws.Range("D73:E73").NumberFormat = "#,##0"
$ws.Range("D73:E73").Font.Name = "Verdana"
$ws.Range("D73:E73").Font.Size = 12
$ws.Range("D73:E73").Font.Bold = $false
$ws.Range("D73:E73").HorizontalAlignment = -4152
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

$ws.Range("D74:E74").NumberFormat = "#,##0"
$ws.Range("D74:E74").Font.Name = "Verdana"
$ws.Range("D74:E74").Font.Size = 12
$ws.Range("D74:E74").Font.Bold = $false
$ws.Range("D74:E74").HorizontalAlignment = -4152
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

$ws.Range("D75:E75").NumberFormat = "#,##0"
$ws.Range("D75:E75").Font.Name = "Verdana"
$ws.Range("D75:E75").Font.Size = 12
$ws.Range("D75:E75").Font.Bold = $false
$ws.Range("D75:E75").HorizontalAlignment = -4152
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

$ws.Range("D76:E76").NumberFormat = "#,##0"
$ws.Range("D76:E76").Font.Name = "Verdana"
$ws.Range("D76:E76").Font.Size = 12
$ws.Range("D76:E76").Font.Bold = $false
$ws.Range("D76:E76").HorizontalAlignment = -4152
$ws.Range("D76").Value = 3253900
$ws.Range("E76").Value = 3340800

$ws.Range("D77:E77").NumberFormat = "#,##0"
$ws.Range("D77:E77").Font.Name = "Verdana"
$ws.Range("D77:E77").Font.Size = 12
$ws.Range("D77:E77").Font.Bold = $false
$ws.Range("D77:E77").HorizontalAlignment = -4152
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

$ws.Range("D80:E80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D80:E80").Font.Name = "Verdana"
$ws.Range("D80:E80").Font.Size = 12
$ws.Range("D80:E80").Font.Bold = $true
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373

$ws.Range("D81:E81").NumberFormat = "#,##0"
$ws.Range("D81:E81").Font.Name = "Verdana"
$ws.Range("D81:E81").Font.Size = 12
$ws.Range("D81:E81").Font.Bold = $false
$ws.Range("D81:E81").HorizontalAlignment = -4152
$ws.Range("D81").Value = 712800
$ws.Range("E81").Value = 556300

$ws.Range("D83:E83").NumberFormat = "#,##0"
$ws.Range("D83:E83").Font.Name = "Verdana"
$ws.Range("D83:E83").Font.Size = 12
$ws.Range("D83:E83").Font.Bold = $false
$ws.Range("D83:E83").HorizontalAlignment = -4152
$ws.Range("D83").Value = 340100
$ws.Range("E83").Value = 335600

$ws.Range("D84:E84").NumberFormat = "#,##0"
$ws.Range("D84:E84").Font.Name = "Verdana"
$ws.Range("D84:E84").Font.Size = 12
$ws.Range("D84:E84").Font.Bold = $false
$ws.Range("D84:E84").HorizontalAlignment = -4152
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

$ws.Range("D85:E85").NumberFormat = "#,##0"
$ws.Range("D85:E85").Font.Name = "Verdana"
$ws.Range("D85:E85").Font.Size = 12
$ws.Range("D85:E85").Font.Bold = $false
$ws.Range("D85:E85").HorizontalAlignment = -4152
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

$ws.Range("D86:E86").NumberFormat = "#,##0"
$ws.Range("D86:E86").Font.Name = "Verdana"
$ws.Range("D86:E86").Font.Size = 12
$ws.Range("D86:E86").Font.Bold = $false
$ws.Range("D86:E86").HorizontalAlignment = -4152
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

$ws.Range("D87:E87").NumberFormat = "#,##0"
$ws.Range("D87:E87").Font.Name = "Verdana"
$ws.Range("D87:E87").Font.Size = 12
$ws.Range("D87:E87").Font.Bold = $false
$ws.Range("D87:E87").HorizontalAlignment = -4152
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

$ws.Range("D88:E88").NumberFormat = "#,##0"
$ws.Range("D88:E88").Font.Name = "Verdana"
$ws.Range("D88:E88").Font.Size = 12
$ws.Range("D88:E88").Font.Bold = $false
$ws.Range("D88:E88").HorizontalAlignment = -4152
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

$ws.Range("D89:E89").NumberFormat = "#,##0"
$ws.Range("D89:E89").Font.Name = "Verdana"
$ws.Range("D89:E89").Font.Size = 12
$ws.Range("D89:E89").Font.Bold = $false
$ws.Range("D89:E89").HorizontalAlignment = -4152
$ws.Range("D89").Value = 999600
$ws.Range("E89").Value = 914400

$ws.Range("D91:E91").NumberFormat = "#,##0"
$ws.Range("D91:E91").Font.Name = "Verdana"
$ws.Range("D91:E91").Font.Size = 12
$ws.Range("D91:E91").Font.Bold = $false
$ws.Range("D91:E91").HorizontalAlignment = -4152
$ws.Range("D91").Value = -196000
$ws.Range("E91").Value = -251300

$ws.Range("D92:E92").NumberFormat = "#,##0"
$ws.Range("D92:E92").Font.Name = "Verdana"
$ws.Range("D92:E92").Font.Size = 12
$ws.Range("D92:E92").Font.Bold = $false
$ws.Range("D92:E92").HorizontalAlignment = -4152
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

$ws.Range("D93:E93").NumberFormat = "#,##0"
$ws.Range("D93:E93").Font.Name = "Verdana"
$ws.Range("D93:E93").Font.Size = 12
$ws.Range("D93:E93").Font.Bold = $false
$ws.Range("D93:E93").HorizontalAlignment = -4152
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

$ws.Range("D94:E94").NumberFormat = "#,##0"
$ws.Range("D94:E94").Font.Name = "Verdana"
$ws.Range("D94:E94").Font.Size = 12
$ws.Range("D94:E94").Font.Bold = $false
$ws.Range("D94:E94").HorizontalAlignment = -4152
$ws.Range("D94").Value = -27100
$ws.Range("E94").Value = -190700

$ws.Range("D96:E96").NumberFormat = "#,##0"
$ws.Range("D96:E96").Font.Name = "Verdana"
$ws.Range("D96:E96").Font.Size = 12
$ws.Range("D96:E96").Font.Bold = $false
$ws.Range("D96:E96").HorizontalAlignment = -4152
$ws.Range("D96").Value = -619500
$ws.Range("E96").Value = -619400

$ws.Range("D97:E97").NumberFormat = "#,##0"
$ws.Range("D97:E97").Font.Name = "Verdana"
$ws.Range("D97:E97").Font.Size = 12
$ws.Range("D97:E97").Font.Bold = $false
$ws.Range("D97:E97").HorizontalAlignment = -4152
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

$ws.Range("D98:E98").NumberFormat = "#,##0"
$ws.Range("D98:E98").Font.Name = "Verdana"
$ws.Range("D98:E98").Font.Size = 12
$ws.Range("D98:E98").Font.Bold = $false
$ws.Range("D98:E98").HorizontalAlignment = -4152
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

$ws.Range("D99:E99").NumberFormat = "#,##0"
$ws.Range("D99:E99").Font.Name = "Verdana"
$ws.Range("D99:E99").Font.Size = 12
$ws.Range("D99:E99").Font.Bold = $false
$ws.Range("D99:E99").HorizontalAlignment = -4152
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

$ws.Range("D100:E100").NumberFormat = "#,##0"
$ws.Range("D100:E100").Font.Name = "Verdana"
$ws.Range("D100:E100").Font.Size = 12
$ws.Range("D100:E100").Font.Bold = $false
$ws.Range("D100:E100").HorizontalAlignment = -4152
$ws.Range("D100").Value = -1154000
$ws.Range("E100").Value = -742200

$ws.Range("D101:E101").NumberFormat = "#,##0"
$ws.Range("D101:E101").Font.Name = "Verdana"
$ws.Range("D101:E101").Font.Size = 12
$ws.Range("D101:E101").Font.Bold = $false
$ws.Range("D101:E101").HorizontalAlignment = -4152
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0

$ws.Range("D102:E102").NumberFormat = "#,##0"
$ws.Range("D102:E102").Font.Name = "Verdana"
$ws.Range("D102:E102").Font.Size = 12
$ws.Range("D102:E102").Font.Bold = $false
$ws.Range("D102:E102").HorizontalAlignment = -4152
$ws.Range("D102").Value = -181400
$ws.Range("E102").Value = -18500

